# Generate Report for Handoff
#
# The nightly localization-status report re-ran and produced fresh
# "Latest Handoff" / "Latest HO Xliff Generate Date" timestamps for the
# e133d234-96c1-4c69-ae7b-10d99374254d.md row (row 7 of the data, i.e.
# sheet row 6 once you count the header). Update the three cells that
# actually changed, leaving everything else untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-22 00:53:17"

# --- zh-cn sheet: "Latest Handoff Datetime" column (H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-22 00:53:13"

# --- de-de sheet: "Latest Handoff Datetime" column (H) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-22 00:53:17"
